# Fixes for merge with MVP-branche / MVPcase-branche
$wb = $excel.ActiveWorkbook

# storageAssets is the 4th sheet (active tab)
$ws = $wb.Worksheets.Item("storageAssets")
$ws.Activate()

# Swap / update rows 14 and 15 so that the EHGV entry ends up in row 15
# and the (renamed) Grid_battery entry ends up in row 14, with updated
# values for the battery row.

# Capture current EHGV row (row 14) values before overwriting
# (use Value2 for reads - Value getter is not reliable in this runtime)
$ehgvName = $ws.Range("B14").Value2
$ehgvType = $ws.Range("D14").Value2

# Row 14 becomes the (renamed) Grid_battery row
$ws.Range("B14").Value = "Grid_battery"
$ws.Range("D14").Value = "STORAGE_ELECTRIC"
$ws.Range("E14").Value = 1000
$ws.Range("L14").Value = 1000
$ws.Range("M14").NumberFormat = "0.00E+00"

# Row 15 becomes the EHGV row
$ws.Range("B15").Value = $ehgvName
$ws.Range("D15").Value = $ehgvType
$ws.Range("E15").Value = 110
$ws.Range("M15").Style = "Normal"

# Update selection to match the recorded view state
$ws.Range("F26").Select()
